$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.236.57"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.898.49"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.08"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4631"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3923"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07873"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9879"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.76"
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("D12").Value = "1.921.18"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.063"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.730"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06984"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.34"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009957"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.03"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "29.255.65"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.300"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.08"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.100"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.84"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.43"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.027"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "118.52"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.885"
$ws.Range("E29").Value = "  -5.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09352"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9038"
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.241"
$ws.Range("E32").Value = "  -2.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.322"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.216"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.185"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05772"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02087"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.712"
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5694"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1782"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.690"
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.88"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5349"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.172"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07026"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.848"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.01"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.056"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.18"
$ws.Range("E51").Value = "  -0.43%  "
